$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like $needle) {
            return $doc.Paragraphs($i)
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Critério de 3 pontos pro terceiro lugar..." paragraph:
#    highlight yellow -> green (paragraph mark + run)
# ---------------------------------------------------------------------
$p1 = Get-ParagraphByText $d "*terceiro lugar, só a partir da rodada 22*"
$p1.Range.Font.HighlightColorIndex = 4

# ---------------------------------------------------------------------
# 2) "Na rodada 38 considerar tbm o critério do saldo de 7 gols." paragraph:
#    highlight yellow -> green (paragraph mark + runs), and merge the
#    three split runs (around the proofErr-wrapped "tbm") into a single run.
# ---------------------------------------------------------------------
$p2 = Get-ParagraphByText $d "*Na rodada 38 considerar*saldo de 7 gols*"
$full2 = $p2.Range
$content2 = $d.Range($full2.Start, $full2.End - 1)

# remember the language formatting carried by the existing text
$lang2 = $d.Range($content2.Start, $content2.Start + 1).LanguageID

# force a genuine content change so the existing run split collapses into one run
$content2.Text = ""
$content2.Text = "Na rodada 38 considerar tbm o critério do saldo de 7 gols."

# re-acquire the paragraph range (content length changed) and restyle it,
# including the paragraph mark, so w:pPr/w:rPr also turns green
$p2b = Get-ParagraphByText $d "Na rodada 38 considerar tbm o critério do saldo de 7 gols.*"
$full2b = $p2b.Range
$full2b.Font.HighlightColorIndex = 4
$full2b.LanguageID = $lang2

# ---------------------------------------------------------------------
# 3) "18) Jogo do oitavo colocado." -> "18) Jogo do sétimo colocado."
#    split into three runs: "18) Jogo do ", "sétimo", " colocado."
# ---------------------------------------------------------------------
$p3 = Get-ParagraphByText $d "*Jogo do oitavo colocado*"
$r3 = $p3.Range
$r3.MoveEnd(1, -1) | Out-Null

$find3 = $r3.Find
$find3.ClearFormatting()
$find3.Text = "oitavo"
$find3.Execute() | Out-Null

$target3 = $find3.Parent
$target3.Text = "sétimo"
# toggling a character attribute keeps this run distinct from its
# (identically formatted) neighbours instead of being silently re-merged
$target3.Bold = 1
$target3.Bold = 0
